# Insert a new data row before row 88 (shifts existing rows 88-173 down to 89-174)
# and populate it with the new "Haba" price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(88).Insert()

$ws.Cells.Item(88, 1).Value  = 3
$ws.Cells.Item(88, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(88, 3).Value  = "Coquimbo"
$ws.Cells.Item(88, 4).Value  = 44778
$ws.Cells.Item(88, 5).Value  = 5
$ws.Cells.Item(88, 6).Value  = 100112026
$ws.Cells.Item(88, 7).Value  = "Haba"
$ws.Cells.Item(88, 8).Value  = "Sin especificar"
$ws.Cells.Item(88, 9).Value  = "Primera"
$ws.Cells.Item(88, 10).Value = 85
$ws.Cells.Item(88, 11).Value = 14000
$ws.Cells.Item(88, 12).Value = 15000
$ws.Cells.Item(88, 13).Value = 14529
$ws.Cells.Item(88, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(88, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(88, 16).Value = 581
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
